# Generate Report for Handoff
# Rename the e2e test markdown file from its old UUID to a new UUID, update
# the associated xliff handoff file names, and bump the "Latest ... Datetime"
# timestamps to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$oldUuid = "2387e6c4-5af5-4cf2-b219-3d1ee8eff66b"
$newUuid = "a17ed91f-f284-4500-9997-cfe95faabba2"

$oldHash = "f9860bc7bf03970f9d3fa77206415ce8da729876"
$newHash = "dba094827919bf27910bc3a66c034aff4cda4e3d"

# Same external link target for every hyperlink in the workbook (unchanged
# by this edit - only the displayed text is updated to the new file name).
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/171cbf24b6961fac54ba1584e15282d272001be1/e2e/$oldUuid.md"

# NOTE: this COM host only binds positional parameters reliably for
# user-defined functions, so the helper below takes its arguments
# positionally (ws, cellRef, displayText) rather than by name.
function Update-DisplayHyperlink {
    param($ws, $cellRef, $displayText)
    $ws.Range($cellRef).Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellRef), $linkAddress, "", "", $displayText) | Out-Null
}

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid.md"
Update-DisplayHyperlink $wsOverview "B2" "e2e\$newUuid.md"
$wsOverview.Range("G2").Value = "2016-08-26 10:55:26"

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid.md"
Update-DisplayHyperlink $wsZhCn "A2" "$newUuid.md"
$wsZhCn.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 10:55:22"

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid.md"
Update-DisplayHyperlink $wsDeDe "A2" "$newUuid.md"
$wsDeDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"
# Latest Handback DateTime (H2) on this sheet is unchanged by this edit.

Write-Host "Updated handoff file references and timestamps."
